$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) from "Neurology" to "Session"
$ws.Name = "Session"

# New rows of data to append (text values, matching existing column layout)
$data = @(
    @("190333", "Neurology", "16/12/2025", "10:13:46", "Manual", "emp17.farah.a.youssef@gmail.com"),
    @("191007", "Neurology", "16/12/2025", "10:16:24", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("202051", "Neurology", "16/12/2025", "10:48:43", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("210728", "Neurology", "16/12/2025", "10:52:36", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("212075", "Neurology", "16/12/2025", "11:10:38", "Manual", "emp17.farah.a.youssef@gmail.com")
)

$startRow = 37
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 1) {
            # Column A values are purely numeric-looking (e.g. "190333") and
            # would otherwise be auto-converted to numbers by Excel. Force
            # them to be stored as text, matching the source data, then
            # reset the style so no extra per-cell style index is applied.
            $cell.NumberFormat = "@"
            $cell.Value = $values[$col - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$col - 1]
        }
    }
}
